$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 339.83334
$ws.Range("I4").Value = 339.83334
$ws.Range("K4").Value = 339.83334
$ws.Range("M4").Value = -225.83334
$ws.Range("H6").Value = 70
$ws.Range("I6").Value = 70
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 210
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -98
$ws.Range("N6").ClearContents()
$ws.Range("H8").Value = 32.5
$ws.Range("I8").Value = 15
$ws.Range("K8").Value = 45
$ws.Range("M8").Value = 94
$ws.Range("H17").Value = 1519.6296
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1519.6296
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4558.8888
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4894.8888
$ws.Range("H19").Value = 921.55554
$ws.Range("I19").Value = 1208.2
$ws.Range("J19").Value = 563.25
$ws.Range("K19").Value = 1208.2
$ws.Range("L19").Value = 563.25
$ws.Range("M19").Value = -1033.2
$ws.Range("N19").Value = -913.25
$ws.Range("H55").Value = 427.1111
$ws.Range("I55").Value = 285.42856
$ws.Range("J55").Value = 517.2727
$ws.Range("K55").Value = 285.42856
$ws.Range("L55").Value = 517.2727
$ws.Range("M55").Value = -71.42856
$ws.Range("N55").Value = -945.2727
$ws.Range("H62").Value = 4931.6665
$ws.Range("I62").Value = 4897.5
$ws.Range("K62").Value = 4897.5
$ws.Range("M62").Value = -4273.5
$ws.Range("H65").Value = 4931.6665
$ws.Range("I65").Value = 4897.5
$ws.Range("K65").Value = 24487.5
$ws.Range("M65").Value = -21367.5
$ws.Range("H80").Value = 1816.6
$ws.Range("I80").Value = 1394.6666
$ws.Range("K80").Value = 4183.9998
$ws.Range("M80").Value = -3185.9998
$ws.Range("H83").Value = 1816.6
$ws.Range("I83").Value = 1394.6666
$ws.Range("K83").Value = 12551.9994
$ws.Range("M83").Value = -7559.999400000001
$ws.Range("H88").Value = 1774.7
$ws.Range("I88").Value = 1193
$ws.Range("K88").Value = 1193
$ws.Range("M88").Value = -787
$ws.Range("H91").Value = 1774.7
$ws.Range("I91").Value = 1193
$ws.Range("K91").Value = 1193
$ws.Range("M91").Value = 211
$ws.Range("H99").Value = 1250
$ws.Range("I99").Value = 1250
$ws.Range("K99").Value = 3750
$ws.Range("M99").Value = -2252
$ws.Range("H111").Value = 3361.6667
$ws.Range("J111").Value = 3209.5
$ws.Range("L111").Value = 9628.5
$ws.Range("N111").Value = -15762.5
$ws.Range("H116").Value = 3595.8572
$ws.Range("J116").Value = 3445.3333
$ws.Range("L116").Value = 3445.3333
$ws.Range("N116").Value = -10329.3333
$ws.Range("H121").Value = 3577.5
$ws.Range("J121").Value = 3577.5
$ws.Range("L121").Value = 10732.5
$ws.Range("N121").Value = -14226.5
$ws.Range("H132").Value = 1578.2593
$ws.Range("I132").Value = 1492.2084
$ws.Range("K132").Value = 4476.6252
$ws.Range("M132").Value = -1946.6252
$ws.Range("H135").Value = 1734.75
$ws.Range("I135").Value = 1508.5
$ws.Range("K135").Value = 13576.5
$ws.Range("M135").Value = -11041.5
$ws.Range("H137").Value = 6586.1665
$ws.Range("I137").Value = 5938
$ws.Range("J137").Value = 7049.143
$ws.Range("K137").Value = 17814
$ws.Range("L137").Value = 21147.429
$ws.Range("M137").Value = -15264
$ws.Range("N137").Value = -26247.429
$ws.Range("H138").Value = 2392.923
$ws.Range("I138").Value = 2758.2
$ws.Range("K138").Value = 8274.599999999999
$ws.Range("M138").Value = -3134.599999999999
$ws.Range("H141").Value = 1738.6
$ws.Range("I141").Value = 1738.6
$ws.Range("K141").Value = 5215.799999999999
$ws.Range("M141").Value = -35.79999999999927

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 10000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H26").Value = 1500
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 1500
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 1500
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -2160
$ws.Range("H32").Value = 8352.75
$ws.Range("I32").Value = 5175.6113
$ws.Range("J32").Value = 17884.166
$ws.Range("K32").Value = 5175.6113
$ws.Range("L32").Value = 17884.166
$ws.Range("M32").Value = -4888.6113
$ws.Range("N32").Value = -18458.166
$ws.Range("H41").Value = 1018.3333
$ws.Range("I41").Value = 1018.3333
$ws.Range("K41").Value = 1018.3333
$ws.Range("M41").Value = -604.3333
$ws.Range("H61").Value = 4513.75
$ws.Range("I61").Value = 4134
$ws.Range("K61").Value = 4134
$ws.Range("M61").Value = -3922
$ws.Range("H74").Value = 1996.1428
$ws.Range("I74").Value = 1596
$ws.Range("K74").Value = 1596
$ws.Range("M74").Value = -722
$ws.Range("H77").Value = 1996.1428
$ws.Range("I77").Value = 1596
$ws.Range("K77").Value = 7980
$ws.Range("M77").Value = -3612
$ws.Range("H132").Value = 5218.36
$ws.Range("I132").Value = 4319.476
$ws.Range("K132").Value = 12958.428
$ws.Range("M132").Value = -10428.428
$ws.Range("H136").Value = 4513.75
$ws.Range("I136").Value = 4134
$ws.Range("K136").Value = 12402
$ws.Range("M136").Value = -9852

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2736.182
$ws.Range("I86").Value = 2721.75
$ws.Range("J86").Value = 2774.6667
$ws.Range("K86").Value = 2721.75
$ws.Range("L86").Value = 2774.6667
$ws.Range("M86").Value = -1598.75
$ws.Range("N86").Value = -5020.6667
$ws.Range("H89").Value = 2736.182
$ws.Range("I89").Value = 2721.75
$ws.Range("J89").Value = 2774.6667
$ws.Range("K89").Value = 13608.75
$ws.Range("L89").Value = 13873.3335
$ws.Range("M89").Value = -7992.75
$ws.Range("N89").Value = -25105.3335
$ws.Range("H105").Value = 4864.65
$ws.Range("I105").Value = 4083.5
$ws.Range("K105").Value = 4083.5
$ws.Range("M105").Value = -2336.5
$ws.Range("H134").Value = 4588
$ws.Range("I134").Value = 2575.8667
$ws.Range("K134").Value = 7727.6001
$ws.Range("M134").Value = -5192.6001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7284.048
$ws.Range("I31").Value = 3111
$ws.Range("J31").Value = 8953.267
$ws.Range("K31").Value = 3111
$ws.Range("L31").Value = 8953.267
$ws.Range("M31").Value = -2816
$ws.Range("N31").Value = -9543.267
$ws.Range("H34").Value = 7284.048
$ws.Range("I34").Value = 3111
$ws.Range("J34").Value = 8953.267
$ws.Range("K34").Value = 3111
$ws.Range("L34").Value = 8953.267
$ws.Range("M34").Value = -2909
$ws.Range("N34").Value = -9357.267
$ws.Range("H58").Value = 6572.0557
$ws.Range("I58").Value = 6456.3125
$ws.Range("K58").Value = 6456.3125
$ws.Range("M58").Value = -6253.3125
$ws.Range("H105").Value = 3793.2354
$ws.Range("I105").Value = 3359
$ws.Range("J105").Value = 4413.5713
$ws.Range("K105").Value = 3359
$ws.Range("L105").Value = 4413.5713
$ws.Range("M105").Value = -1612
$ws.Range("N105").Value = -7907.5713
$ws.Range("H107").Value = 1118.4
$ws.Range("I107").Value = 1248.75
$ws.Range("K107").Value = 1248.75
$ws.Range("M107").Value = 671.25
$ws.Range("H134").Value = 5919.294
$ws.Range("I134").Value = 5769.5483
$ws.Range("J134").Value = 7466.6665
$ws.Range("K134").Value = 17308.6449
$ws.Range("L134").Value = 22399.9995
$ws.Range("M134").Value = -14773.6449
$ws.Range("N134").Value = -27469.9995
$ws.Range("H135").Value = 99999
$ws.Range("J135").Value = 99998
$ws.Range("L135").Value = 99998
$ws.Range("N135").Value = -110138
$ws.Range("H136").Value = 6572.0557
$ws.Range("I136").Value = 6456.3125
$ws.Range("K136").Value = 19368.9375
$ws.Range("M136").Value = -16818.9375
$ws.Range("H140").Value = 99994.5
$ws.Range("J140").Value = 99994.5
$ws.Range("L140").Value = 99994.5
$ws.Range("N140").Value = -110354.5
$ws.Range("H141").Value = 314324.66
$ws.Range("J141").Value = 314324.66
$ws.Range("L141").Value = 314324.66
$ws.Range("N141").Value = -324684.66

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 250708.75
$ws.Range("J32").Value = 945
$ws.Range("L32").Value = 2835
$ws.Range("N32").Value = -3401
$ws.Range("H37").Value = 97897
$ws.Range("J37").Value = 97897
$ws.Range("L37").Value = 293691
$ws.Range("N37").Value = -293915
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H61").Value = 59.4
$ws.Range("I61").Value = 95
$ws.Range("J61").Value = 50.5
$ws.Range("K61").Value = 285
$ws.Range("L61").Value = 151.5
$ws.Range("M61").Value = -70
$ws.Range("N61").Value = -581.5
$ws.Range("H68").Value = 1299.3334
$ws.Range("I68").Value = 949
$ws.Range("K68").Value = 2847
$ws.Range("M68").Value = -2036
$ws.Range("H71").Value = 1299.3334
$ws.Range("I71").Value = 949
$ws.Range("K71").Value = 8541
$ws.Range("M71").Value = -4485
$ws.Range("H80").Value = 15484.143
$ws.Range("J80").Value = 1497.5
$ws.Range("L80").Value = 4492.5
$ws.Range("N80").Value = -6364.5
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H83").Value = 15484.143
$ws.Range("J83").Value = 1497.5
$ws.Range("L83").Value = 13477.5
$ws.Range("N83").Value = -22837.5
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H112").Value = 945
$ws.Range("I112").Value = 945
$ws.Range("K112").Value = 2835
$ws.Range("M112").Value = -1727
$ws.Range("H116").Value = 2405.2727
$ws.Range("J116").Value = 2499.889
$ws.Range("L116").Value = 7499.667
$ws.Range("N116").Value = -14383.667
$ws.Range("H131").Value = 2247.0625
$ws.Range("I131").Value = 1056.625
$ws.Range("J131").Value = 3437.5
$ws.Range("K131").Value = 3169.875
$ws.Range("L131").Value = 10312.5
$ws.Range("M131").Value = 1870.125
$ws.Range("N131").Value = -20392.5
$ws.Range("H137").Value = 2067.7778
$ws.Range("I137").Value = 2067.7778
$ws.Range("K137").Value = 6203.3334
$ws.Range("M137").Value = -1103.3334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5766250.5
$ws.Range("I11").Value = 4345000
$ws.Range("J11").Value = 6619001
$ws.Range("K11").Value = 4345000
$ws.Range("L11").Value = 6619001
$ws.Range("M11").Value = -4344861
$ws.Range("N11").Value = -6619279
$ws.Range("H14").Value = 1670002.6
$ws.Range("I14").Value = 1670002.6
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1670002.6
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1669834.6
$ws.Range("N14").ClearContents()
$ws.Range("H80").Value = 9557.977000000001
$ws.Range("I80").Value = 5787.222
$ws.Range("J80").Value = 12386.042
$ws.Range("K80").Value = 5787.222
$ws.Range("L80").Value = 12386.042
$ws.Range("M80").Value = -4789.222
$ws.Range("N80").Value = -14382.042
$ws.Range("H83").Value = 9557.977000000001
$ws.Range("I83").Value = 5787.222
$ws.Range("J83").Value = 12386.042
$ws.Range("K83").Value = 28936.11
$ws.Range("L83").Value = 61930.21
$ws.Range("M83").Value = -23944.11
$ws.Range("N83").Value = -71914.20999999999
$ws.Range("H113").Value = 1707.1111
$ws.Range("I113").Value = 1487.9286
$ws.Range("K113").Value = 1487.9286
$ws.Range("M113").Value = 682.0714
$ws.Range("H132").Value = 727
$ws.Range("I132").Value = 665.6
$ws.Range("K132").Value = 1996.8
$ws.Range("M132").Value = 533.1999999999998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7142.4
$ws.Range("I7").Value = 7748.2856
$ws.Range("J7").Value = 6612.25
$ws.Range("K7").Value = 7748.2856
$ws.Range("L7").Value = 6612.25
$ws.Range("M7").Value = -7636.2856
$ws.Range("N7").Value = -6836.25
$ws.Range("H22").Value = 313.33334
$ws.Range("I22").Value = 395
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 395
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = -100
$ws.Range("N22").Value = -740
$ws.Range("H26").Value = 4042.8333
$ws.Range("I26").Value = 3000
$ws.Range("J26").Value = 4877.1
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 4877.1
$ws.Range("M26").Value = -2705
$ws.Range("N26").Value = -5467.1
$ws.Range("H27").Value = 313.33334
$ws.Range("I27").Value = 395
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 395
$ws.Range("L27").Value = 150
$ws.Range("M27").Value = -288
$ws.Range("N27").Value = -364
$ws.Range("H40").Value = 3488.3333
$ws.Range("I40").Value = 3100
$ws.Range("J40").Value = 4847.5
$ws.Range("K40").Value = 3100
$ws.Range("L40").Value = 4847.5
$ws.Range("M40").Value = -2964
$ws.Range("N40").Value = -5119.5
$ws.Range("H46").Value = 1249.7222
$ws.Range("I46").Value = 1956.8
$ws.Range("J46").Value = 977.7692
$ws.Range("K46").Value = 1956.8
$ws.Range("L46").Value = 977.7692
$ws.Range("M46").Value = -1768.8
$ws.Range("N46").Value = -1353.7692
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H107").Value = 2993
$ws.Range("I107").Value = 2993
$ws.Range("K107").Value = 2993
$ws.Range("M107").Value = -1073
$ws.Range("H126").Value = 7142.4
$ws.Range("I126").Value = 7748.2856
$ws.Range("J126").Value = 6612.25
$ws.Range("K126").Value = 23244.8568
$ws.Range("L126").Value = 19836.75
$ws.Range("M126").Value = -20774.8568
$ws.Range("N126").Value = -24776.75
$ws.Range("H132").Value = 2512.4736
$ws.Range("I132").Value = 2365.2593
$ws.Range("J132").Value = 2873.818
$ws.Range("K132").Value = 7095.777900000001
$ws.Range("L132").Value = 8621.454000000002
$ws.Range("M132").Value = -4565.777900000001
$ws.Range("N132").Value = -13681.454
$ws.Range("H136").Value = 5036.9614
$ws.Range("I136").Value = 4122.5
$ws.Range("J136").Value = 5820.7856
$ws.Range("K136").Value = 12367.5
$ws.Range("L136").Value = 17462.3568
$ws.Range("M136").Value = -9817.5
$ws.Range("N136").Value = -22562.3568
$ws.Range("H138").Value = 58000
$ws.Range("J138").Value = 58000
$ws.Range("L138").Value = 58000
$ws.Range("N138").Value = -68280
$ws.Range("H141").Value = 66250
$ws.Range("J141").Value = 66250
$ws.Range("L141").Value = 66250
$ws.Range("N141").Value = -76610

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1999.5
$ws.Range("J6").Value = 1999.5
$ws.Range("L6").Value = 1999.5
$ws.Range("N6").Value = -2229.5
$ws.Range("H14").Value = 55498.832
$ws.Range("I14").Value = 9499.5
$ws.Range("K14").Value = 9499.5
$ws.Range("M14").Value = -9331.5
$ws.Range("H58").Value = 18483.334
$ws.Range("I58").Value = 22500
$ws.Range("K58").Value = 22500
$ws.Range("M58").Value = -22192
$ws.Range("H113").Value = 479.7
$ws.Range("I113").Value = 479.7
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1439.1
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 730.9000000000001
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 123630.555
$ws.Range("I122").Value = 12445.833
$ws.Range("J122").Value = 346000
$ws.Range("K122").Value = 37337.499
$ws.Range("L122").Value = 1038000
$ws.Range("M122").Value = -34887.499
$ws.Range("N122").Value = -1042900
$ws.Range("H132").Value = 3457.1538
$ws.Range("I132").Value = 3540.3635
$ws.Range("K132").Value = 10621.0905
$ws.Range("M132").Value = -8091.0905
$ws.Range("H135").Value = 59123.75
$ws.Range("J135").Value = 59123.75
$ws.Range("L135").Value = 59123.75
$ws.Range("N135").Value = -69263.75
$ws.Range("H136").Value = 3161.4
$ws.Range("I136").Value = 3214
$ws.Range("K136").Value = 9642
$ws.Range("M136").Value = -7092
